# "Generate Report for Handback" — refresh the localization-status report
# after a successful handback that is now in sync with en-US.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ----------------------------------------------------------------------
# Overview sheet: roll the per-locale status up from "Ready for handoff"
# to the new "Handed back" status.
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus   # zh-cn
$overview.Range("F2").Value = $newStatus   # de-de

# Widen the two status columns so the longer message is readable.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ----------------------------------------------------------------------
# zh-cn sheet: update status, bump the handback timestamp, and clear the
# stale "handback not latest" error now that it is in sync.
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-26 08:54:19"
$zhcn.Range("P2").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ----------------------------------------------------------------------
# de-de sheet: same treatment.
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-26 08:54:26"
$dede.Range("P2").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
